function Set-TextValue($range, $value) {
    # Force the cell to remain plain text even when the literal
    # looks numeric (e.g. "569.20") or looks like a multi-dot
    # price string (e.g. "64.249.88"). Resetting the style back to
    # "Normal" afterwards drops the temporary Text number format so
    # we do not leave a stray style index on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range('D2') '64.249.88'
Set-TextValue $ws.Range('E2') '  -3.38%  '

Set-TextValue $ws.Range('D3') '3.156.95'
Set-TextValue $ws.Range('E3') '  -5.04%  '

Set-TextValue $ws.Range('E4') '  +0.24%  '

Set-TextValue $ws.Range('D5') '569.20'
Set-TextValue $ws.Range('E5') '  -3.00%  '

Set-TextValue $ws.Range('D6') '167.57'
Set-TextValue $ws.Range('E6') '  -7.47%  '

Set-TextValue $ws.Range('D7') '0.602'
Set-TextValue $ws.Range('E7') '  -7.85%  '

Set-TextValue $ws.Range('E8') '  +0.00%  '

Set-TextValue $ws.Range('D9') '3.184.88'
Set-TextValue $ws.Range('E9') '  -4.06%  '

Set-TextValue $ws.Range('D10') '0.118'
Set-TextValue $ws.Range('E10') '  -6.66%  '

Set-TextValue $ws.Range('D11') '6.81'
Set-TextValue $ws.Range('E11') '  -0.29%  '

Set-TextValue $ws.Range('D12') '0.387'
Set-TextValue $ws.Range('E12') '  -3.75%  '

Set-TextValue $ws.Range('D13') '3.707.38'
Set-TextValue $ws.Range('E13') '  -4.90%  '

Set-TextValue $ws.Range('E14') '  -1.37%  '

Set-TextValue $ws.Range('D15') '64.376.69'
Set-TextValue $ws.Range('E15') '  -3.12%  '

Set-TextValue $ws.Range('D16') '25.24'
Set-TextValue $ws.Range('E16') '  -5.03%  '

Set-TextValue $ws.Range('B17') 'WrappedEther'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D17') '3.166.69'
Set-TextValue $ws.Range('E17') '  -3.92%  '

Set-TextValue $ws.Range('B18') 'ShibaInu'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D18') '0.0000156'
Set-TextValue $ws.Range('E18') '  -5.01%  '

Set-TextValue $ws.Range('D19') '415.95'
Set-TextValue $ws.Range('E19') '  -2.12%  '

Set-TextValue $ws.Range('D20') '12.81'
Set-TextValue $ws.Range('E20') '  -2.68%  '

Set-TextValue $ws.Range('D21') '5.28'
Set-TextValue $ws.Range('E21') '  -3.99%  '

Set-TextValue $ws.Range('D22') '7.13'
Set-TextValue $ws.Range('E22') '  -3.38%  '

Set-TextValue $ws.Range('D23') '0.999'
Set-TextValue $ws.Range('E23') '  -0.20%  '

Set-TextValue $ws.Range('E24') '  +0.00%  '

Set-TextValue $ws.Range('D25') '69.57'
Set-TextValue $ws.Range('E25') '  -3.24%  '

Set-TextValue $ws.Range('E26') '  -0.67%  '

Set-TextValue $ws.Range('D27') '0.497'
Set-TextValue $ws.Range('E27') '  -3.55%  '

Set-TextValue $ws.Range('D28') '0.0000102'
Set-TextValue $ws.Range('E28') '  -11.21%  '

Set-TextValue $ws.Range('D29') '8.77'
Set-TextValue $ws.Range('E29') '  -3.54%  '

Set-TextValue $ws.Range('D30') '0.998'
Set-TextValue $ws.Range('E30') '  -0.04%  '

Set-TextValue $ws.Range('D31') '1.82'
Set-TextValue $ws.Range('E31') '  -5.27%  '

Set-TextValue $ws.Range('D32') '21.70'
Set-TextValue $ws.Range('E32') '  -3.06%  '

Set-TextValue $ws.Range('E33') '  -0.14%  '

Set-TextValue $ws.Range('D34') '5.03'
Set-TextValue $ws.Range('E34') '  -3.01%  '

Set-TextValue $ws.Range('D35') '6.34'
Set-TextValue $ws.Range('E35') '  -4.17%  '

Set-TextValue $ws.Range('E36') '  -5.91%  '

Set-TextValue $ws.Range('D37') '155.35'
Set-TextValue $ws.Range('E37') '  -3.14%  '

Set-TextValue $ws.Range('E38') '  -5.73%  '

Set-TextValue $ws.Range('D39') '2.694.26'
Set-TextValue $ws.Range('E39') '  -6.08%  '

Set-TextValue $ws.Range('D40') '1.69'
Set-TextValue $ws.Range('E40') '  -6.77%  '

Set-TextValue $ws.Range('D41') '4.19'
Set-TextValue $ws.Range('E41') '  -3.21%  '

Set-TextValue $ws.Range('D42') '23.84'
Set-TextValue $ws.Range('E42') '  -9.71%  '

Set-TextValue $ws.Range('D43') '39.04'
Set-TextValue $ws.Range('E43') '  -1.89%  '

Set-TextValue $ws.Range('E44') '  -5.68%  '

Set-TextValue $ws.Range('D45') '0.0614'
Set-TextValue $ws.Range('E45') '  -7.26%  '

Set-TextValue $ws.Range('D46') '5.46'
Set-TextValue $ws.Range('E46') '  -7.64%  '

Set-TextValue $ws.Range('D47') '0.0261'
Set-TextValue $ws.Range('E47') '  -4.46%  '

Set-TextValue $ws.Range('D48') '287.94'
Set-TextValue $ws.Range('E48') '  -8.22%  '

Set-TextValue $ws.Range('D49') '21.18'
Set-TextValue $ws.Range('E49') '  -8.75%  '

Set-TextValue $ws.Range('E50') '  +0.30%  '

Set-TextValue $ws.Range('D51') '0.0987'
Set-TextValue $ws.Range('E51') '  -5.79%  '
